$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. F8: Open -> Closed -------------------------------------------------
$ws.Range("F8").Value = "Closed"

# --- 2. Row 9: "Move SRS context under project description" --------------
$ws.Range("A2").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A9").Value = "24/2/2020"

$ws.Range("B2").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("B9").Value = "Ali"

$ws.Range("C2").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("C9").Value = "SRS"

$ws.Range("D9").HorizontalAlignment = -4131
$ws.Range("D9").VerticalAlignment = -4160
$ws.Range("D9").Value = "Move SRS context under project description"

$ws.Range("E2").Copy()
$ws.Range("E9").PasteSpecial(-4122)
$ws.Range("E9").ClearContents()

$ws.Range("F2").Copy()
$ws.Range("F9").PasteSpecial(-4122)
$ws.Range("F9").Value = "Open"

# --- 3. Row 10: "SRS context shall not has a req_ID ..." -----------------
$ws.Range("A2").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Value = "24/2/2020"

$ws.Range("B2").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("B10").Value = "Ali"

$ws.Range("C2").Copy()
$ws.Range("C10").PasteSpecial(-4122)
$ws.Range("C10").Value = "SRS"

$ws.Range("D9").Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("D10").WrapText = $true
$ws.Range("D10").Value = "SRS context shall not has a req_ID remove it's presence form `nSRS_013, alos inputs and outputs have no meaning for this requirement"

$ws.Range("E2").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$ws.Range("E10").ClearContents()

$ws.Range("F2").Copy()
$ws.Range("F10").PasteSpecial(-4122)
$ws.Range("F10").Value = "Open"

$ws.Rows.Item(10).RowHeight = 43.2

# --- 4. Row 11: "Req_PO1_DGC_SRS_001_V01: ..." ----------------------------
$ws.Range("A7").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = "28/2/2020"

$ws.Range("B2").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("B11").Value = "Alzahraa"

$ws.Range("C2").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("C11").Value = "SRS"

$ws.Range("D2").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("D11").Value = "Req_PO1_DGC_SRS_001_V01:`nRequirement shall be splitted into 2 requirements, one for accepting input op1, op2 and operator, and another for calculating result"

$ws.Range("F2").Copy()
$ws.Range("F11").PasteSpecial(-4122)
$ws.Range("F11").Value = "Open"

$ws.Rows.Item(11).RowHeight = 43.2

# --- 5. Conditional formatting for the new rows (same scheme as F2/E2 etc.) -
$fc = $ws.Range("F9").FormatConditions
$c = $fc.Add(9, $null, $null, $null, "Open", 0)
$c.Font.Color = 255
$c = $fc.Add(9, $null, $null, $null, "Closed", 0)
$c.Font.Color = 24832
$c.Interior.Color = 13561798
$c = $fc.Add(9, $null, $null, $null, "Open", 0)
$c.Font.Color = 393372

$fc = $ws.Range("E9").FormatConditions
$c = $fc.Add(1, 3, '"Rejected"')
$c.Font.Color = 393372
$c = $fc.Add(1, 3, '"Accepted"')
$c.Font.Color = 24832
$c.Interior.Color = 13561798

$fc = $ws.Range("F10").FormatConditions
$c = $fc.Add(9, $null, $null, $null, "Open", 0)
$c.Font.Color = 255
$c = $fc.Add(9, $null, $null, $null, "Closed", 0)
$c.Font.Color = 24832
$c.Interior.Color = 13561798
$c = $fc.Add(9, $null, $null, $null, "Open", 0)
$c.Font.Color = 393372

$fc = $ws.Range("E10").FormatConditions
$c = $fc.Add(1, 3, '"Rejected"')
$c.Font.Color = 393372
$c = $fc.Add(1, 3, '"Accepted"')
$c.Font.Color = 24832
$c.Interior.Color = 13561798

$fc = $ws.Range("F11").FormatConditions
$c = $fc.Add(9, $null, $null, $null, "Open", 0)
$c.Font.Color = 255
$c = $fc.Add(9, $null, $null, $null, "Closed", 0)
$c.Font.Color = 24832
$c.Interior.Color = 13561798
$c = $fc.Add(9, $null, $null, $null, "Open", 0)
$c.Font.Color = 393372

# --- 6. Extend the data-validation dropdowns to cover the new rows --------
$ws.Range("F2:F8").Validation.Delete()
$ws.Range("F2:F11").Validation.Add(3, 1, 1, '"Open, Closed"')

$ws.Range("E2:E8").Validation.Delete()
$ws.Range("E2:E10").Validation.Add(3, 1, 1, '"Accepted, Rejected"')

# --- 7. Update the saved view / selection ----------------------------------
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("D12").Select()
